$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix typo in Numero_Documento for row 6 (extra trailing digit removed).
$ws.Range("B6").Value = 1058550941

# Insert a new row at position 53. This shifts the existing row 53
# (ODON-0000040 / BRESLI TATIANA LEGARDA / ...) down to row 54 intact.
$ws.Rows.Item(53).Insert()

# Populate the new row 53 with the new appointment record
# (ODON-0000005 / EMANUEL GAON MELLIZO / PREVIRED / Ortodoncia / Efectivo).
$ws.Range("A53").Value = "ODON-0000005"
$ws.Range("B53").Value = 1058550941
$ws.Range("C53").Value = "EMANUEL GAON MELLIZO"
$ws.Range("E53").Value = "PREVIRED"

# Force F53 to stay plain text ("04/02/2026") instead of being
# auto-converted to a date serial value by COM's type inference.
$ws.Range("F53").NumberFormat = "@"
$ws.Range("F53").Value = "04/02/2026"
$ws.Range("F53").Style = "Normal"

$ws.Range("G53").Value = 2026
$ws.Range("H53").Value = "FEBRERO"
$ws.Range("I53").Value = "SEMANA1"
$ws.Range("J53").Value = "LINA ACOSTA"
$ws.Range("K53").Value = "MARLEN MUÑOZ"
$ws.Range("L53").Value = "Ortodoncia"
$ws.Range("N53").Value = "Procedimiento"
$ws.Range("O53").Value = 1
$ws.Range("P53").Value = 1
$ws.Range("Q53").Value = 1
$ws.Range("V53").Value = "Efectivo"
$ws.Range("W53").Value = "LINA ACOSTA"
$ws.Range("AA53").Value = 275000
